# Individual Contribution Metrics - add per-member breakdown for newly
# completed tasks (rows 8, 11, 12, 13) and append the new tasks that were
# delivered (rows 14-19), each with its per-member score breakdown in
# columns E (Rodrigo), F (Jose) and G (Cesar). The summary formulas in
# I3:K3 / I4:K4 recalc automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Existing rows: fill in newly-attributed per-member scores ---------

# Row 8 - Stress and Anxiety test
$ws.Range("E8").Value = 5
$ws.Range("F8").Value = 5

# Row 11 - Research ways to handle the follow-up of the users
$ws.Range("E11").Value = 8
$ws.Range("F11").Value = 8
$ws.Range("G11").Value = 8

# Row 12 - Add ways to contact a pyschologist
$ws.Range("F12").Value = 2

# Row 13 - Documentation of the project (Project logs)
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0

# --- New tasks -----------------------------------------------------------
# Write the B-column task names in the order that makes the new shared
# strings land in the same order as the source workbook (Usability Tests,
# Final interfaces, HCI Article, HCI Poster, HCI Video, Usability Test
# Artifacts).

$ws.Range("B16").Value = "Usability Tests"
$ws.Range("B14").Value = "Final interfaces"
$ws.Range("B17").Value = "HCI Article"
$ws.Range("B18").Value = "HCI Poster"
$ws.Range("B19").Value = "HCI Video"
$ws.Range("B15").Value = "Usability Test Artifacts (script, checklist)"

# Row 14 - Final interfaces
$ws.Range("A14").Value = 12
$ws.Range("C14").Value = 13
$ws.Range("F14").Value = 13

# Row 15 - Usability Test Artifacts (script, checklist)
$ws.Range("A15").Value = 13
$ws.Range("C15").Value = 5
$ws.Range("E15").Value = 5
$ws.Range("G15").Value = 5

# Row 16 - Usability Tests
$ws.Range("A16").Value = 14
$ws.Range("C16").Value = 13
$ws.Range("E16").Value = 13

# Row 17 - HCI Article
$ws.Range("A17").Value = 15
$ws.Range("C17").Value = 8
$ws.Range("F17").Value = 8
$ws.Range("G17").Value = 8

# Row 18 - HCI Poster
$ws.Range("A18").Value = 16
$ws.Range("C18").Value = 3
$ws.Range("E18").Value = 3
$ws.Range("F18").Value = 3

# Row 19 - HCI Video
$ws.Range("A19").Value = 17
$ws.Range("C19").Value = 5
$ws.Range("E19").Value = 5
$ws.Range("F19").Value = 5
$ws.Range("G19").Value = 5

# --- Selection / view state ----------------------------------------------
$ws.Range("K11").Select()
